$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.379.33'
$ws.Range('E2').Value = '  -0.32%  '
$ws.Range('D3').Value = '1.565.56'
$ws.Range('E3').Value = '  -0.21%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = "'210.58"
$ws.Range('E5').Value = '  -0.58%  '
$ws.Range('D6').Value = "'0.489"
$ws.Range('E6').Value = '  -0.40%  '
$ws.Range('D8').Value = "'44.23"
$ws.Range('E8').Value = '  -4.33%  '
$ws.Range('D9').Value = "'23.51"
$ws.Range('E9').Value = '  -2.09%  '
$ws.Range('E10').Value = '  -1.14%  '
$ws.Range('E11').Value = '  -0.66%  '
$ws.Range('D12').Value = "'0.0890"
$ws.Range('E12').Value = '  +0.03%  '
$ws.Range('D13').Value = '1.788.25'
$ws.Range('E13').Value = '  -0.21%  '
$ws.Range('D14').Value = '1.556.04'
$ws.Range('E14').Value = '  -0.83%  '
$ws.Range('B15').Value = 'WrappedBTC'
$ws.Range('C15').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D15').Value = '28.349.57'
$ws.Range('E15').Value = '  -0.43%  '
$ws.Range('B16').Value = 'Polkadot'
$ws.Range('C16').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D16').Value = "'3.66"
$ws.Range('E16').Value = '  -0.58%  '
$ws.Range('E17').Value = '  -1.71%  '
$ws.Range('D18').Value = "'60.39"
$ws.Range('E18').Value = '  -2.98%  '
$ws.Range('D19').Value = "'227.44"
$ws.Range('E19').Value = '  -0.14%  '
$ws.Range('D20').Value = "'7.36"
$ws.Range('E20').Value = '  +0.25%  '
$ws.Range('D21').Value = '0.0₃0677'
$ws.Range('E21').Value = '  -1.88%  '
$ws.Range('E22').Value = '  -0.05%  '
$ws.Range('D23').Value = "'3.93"
$ws.Range('E23').Value = '  +1.31%  '
$ws.Range('D24').Value = "'8.90"
$ws.Range('E24').Value = '  -2.42%  '
$ws.Range('E25').Value = '  -2.89%  '
$ws.Range('D26').Value = "'150.23"
$ws.Range('E26').Value = '  -0.43%  '
$ws.Range('D27').Value = "'14.86"
$ws.Range('E27').Value = '  -0.75%  '
$ws.Range('D28').Value = "'0.103"
$ws.Range('E28').Value = '  -0.03%  '
$ws.Range('B29').Value = 'BinanceUSD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D29').Value = "'1.00"
$ws.Range('E29').Value = '  -0.09%  '
$ws.Range('B30').Value = 'Cosmos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D30').Value = "'6.28"
$ws.Range('E30').Value = '  -2.66%  '
$ws.Range('D31').Value = "'0.0473"
$ws.Range('E31').Value = '  +0.44%  '
$ws.Range('D32').Value = "'1.07"
$ws.Range('E32').Value = '  -3.59%  '
$ws.Range('E33').Value = '  -0.74%  '
$ws.Range('E34').Value = '  -0.22%  '
$ws.Range('D35').Value = '1.381.59'
$ws.Range('E35').Value = '  -0.82%  '
$ws.Range('D36').Value = "'1.06"
$ws.Range('E36').Value = '  +2.10%  '
$ws.Range('E37').Value = '  -3.91%  '
$ws.Range('E38').Value = '  -0.58%  '
$ws.Range('E39').Value = '  +2.39%  '
$ws.Range('E40').Value = '  -2.25%  '
$ws.Range('B41').Value = 'ImmutableX'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D41').Value = "'0.517"
$ws.Range('E41').Value = '  -4.03%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').Value = "'1.93"
$ws.Range('E42').Value = '  +2.61%  '
$ws.Range('E43').Value = '  -0.09%  '
$ws.Range('D44').Value = "'0.783"
$ws.Range('E44').Value = '  -0.40%  '
$ws.Range('D45').Value = "'0.0465"
$ws.Range('E45').Value = '  -2.45%  '
$ws.Range('D46').Value = "'5.35"
$ws.Range('E46').Value = '  -2.46%  '
$ws.Range('D47').Value = "'62.11"
$ws.Range('E47').Value = '  -1.11%  '
$ws.Range('D48').Value = "'0.915"
$ws.Range('E48').Value = '  -5.99%  '
$ws.Range('D49').Value = '1.701.13'
$ws.Range('E49').Value = '  -0.18%  '
$ws.Range('D50').Value = "'85.32"
$ws.Range('E50').Value = '  -0.69%  '
$ws.Range('E51').Value = '  -1.42%  '
